$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.552.91'
$ws.Range('D3').Value = '1.913.10'
$ws.Range('E3').Value = '  +5.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.24'
$ws.Range('E5').Value = '  +1.41%  '
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('E7').Value = '  +4.98%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3970'
$ws.Range('E8').Value = '  +1.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09694'
$ws.Range('E9').Value = '  -1.61%  '
$ws.Range('E10').Value = '  +4.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.97'
$ws.Range('E11').Value = '  +2.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.542'
$ws.Range('E12').Value = '  +1.68%  '
$ws.Range('E13').Value = '  +3.05%  '
$ws.Range('D14').Value = '1.909.75'
$ws.Range('E14').Value = '  +5.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.565'
$ws.Range('E15').Value = '  +3.91%  '
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.84'
$ws.Range('E17').Value = '  +2.71%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001139'
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06652'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.18'
$ws.Range('E20').Value = '  +5.66%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.342'
$ws.Range('E22').Value = '  +6.47%  '
$ws.Range('D23').Value = '28.643.34'
$ws.Range('E23').Value = '  +1.93%  '
$ws.Range('E24').Value = '  +2.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.306'
$ws.Range('E25').Value = '  +2.65%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.398'
$ws.Range('E26').Value = '  -0.70%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.691'
$ws.Range('E27').Value = '  +11.57%  '
$ws.Range('B28').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C28').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D28').Value = '2.129.76'
$ws.Range('E28').Value = '  +4.96%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '21.26'
$ws.Range('E29').Value = '  +2.89%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '158.79'
$ws.Range('E30').Value = '  -0.27%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '128.93'
$ws.Range('E31').Value = '  +1.69%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.110'
$ws.Range('E32').Value = '  +7.41%  '
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.1086'
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.756'
$ws.Range('E34').Value = '  +3.39%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.637'
$ws.Range('E35').Value = '  +0.53%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.918'
$ws.Range('E36').Value = '  +11.33%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06762'
$ws.Range('E37').Value = '  +1.08%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02440'
$ws.Range('E38').Value = '  +4.34%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.266'
$ws.Range('E39').Value = '  +7.69%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2229'
$ws.Range('E40').Value = '  +4.09%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.84'
$ws.Range('E41').Value = '  +4.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6479'
$ws.Range('E42').Value = '  +4.55%  '
$ws.Range('B43').Value = 'InternetComputer(DFINITY)'
$ws.Range('C43').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.088'
$ws.Range('E43').Value = '  +2.58%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.191'
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.58'
$ws.Range('E46').Value = '  +2.69%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.6103'
$ws.Range('E47').Value = '  +3.26%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.763'
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('B49').Value = 'WEMIXTOKEN'
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.287'
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('B50').Value = 'NEARProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.033'
$ws.Range('E50').Value = '  +4.72%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '125.24'
$ws.Range('E51').Value = '  +0.64%  '
